# "Generate Report for Archive"
#
# The localization status report is regenerated: the outstanding item that
# was "Ready for handoff" has moved on to "In Translation", and the
# (generator-computed) widths of the Status / zh-cn / de-de columns shrink
# to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" --------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# --- Column widths re-fitted to the shorter text ------------------------
# (matches the narrower width now needed for "In Translation")
$overview.Columns.Item(5).ColumnWidth = 12.5   # E: zh-cn
$overview.Columns.Item(6).ColumnWidth = 12.5   # F: de-de
$zhcn.Columns.Item(3).ColumnWidth     = 12.5   # C: Status
$dede.Columns.Item(3).ColumnWidth     = 12.5   # C: Status
